# Fill in column D (F2_Workout_Week_Score) with scores for each team,
# matching the same red-font / thin-border style already used in column C,
# and recolor the last four rows of column C (C10:C13) to the accent5 theme color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to write into D2:D13
$values = @(150, 100, 100, 100, 110, 100, 100, 100, 130, 100, 100, 100)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("D$row")
    $cell.Value = $values[$i]
    # Match the style already used in column C for these rows (red font, thin border)
    $cell.Font.Color = 255
}

# Recolor C10:C13 to the new accent5 theme color font
$ws.Range("C10:C13").Font.ThemeColor = 9

# Update the active selection to D2:D13 (active cell D2)
$ws.Range("D2:D13").Select()
